# daily auto push: 2026-02-28 22:37 UTC
#
# The log table on Sheet1 gains one new entry. A new data row is inserted
# right after the existing "2026/03/01" row (the current row 881), pushing
# every row below it down by one (old 881 -> new 882, ..., old 922 -> new
# 923). The new row repeats the same date ("2026/03/01") and weekday
# ("日") as its neighbour above, with its own time-slot/ranking values
# (4 / 37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 881, shifting rows 881:922 down to 882:923.
$ws.Rows.Item(881).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/12/29"), not
# real Excel date serials. Force text formatting before writing the value
# so "2026/03/01" isn't auto-converted into a date, then drop the leftover
# number-format override so the cell's style matches its plain neighbours.
$newRow = $ws.Rows.Item(881)
$newRow.Cells.Item(1, 1).NumberFormat = "@"
$newRow.Cells.Item(1, 1).Value = "2026/03/01"
$newRow.Cells.Item(1, 1).ClearFormats()

$newRow.Cells.Item(1, 2).Value = "日"
$newRow.Cells.Item(1, 3).Value = 4
$newRow.Cells.Item(1, 4).Value = 37
